$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new column headers I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style (bold font, border, centered alignment) from the existing
# header cell H1 onto the two new header cells using a formats-only paste,
# so the underlying shared style index is reused rather than duplicated.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows 2-43: new columns I ("I0") and J ("IF") ---
$iValues = @(5,8,8,7,8,8,10,10,9,8,9,8,8,9,8,9,8,6,9,10,7,8,8,8,7,7,7,7,8,8,8,7,8,3,8,7,9,8,7,8,7,7)
$jValues = @(7,8,8,7,8,8,10,10,9,8,9,8,8,9,8,9,8,7,9,10,7,8,8,8,7,7,8,7,8,8,8,7,8,3,8,8,9,8,8,8,7,8)

for ($r = 2; $r -le 43; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}

Write-Host "done"
